# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet and the
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" values for
# the 1c450ed0-f87e-43b5-96ad-5ccaa6b0fded.md row on the zh-cn and de-de
# per-language report sheets, reflecting a new handback report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
$wsZhCn.Range("H2").Value = "2016-08-24 12:50:50"
$wsZhCn.Range("K2").Value = "2016-08-24 12:51:20"

# de-de: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
$wsDeDe.Range("H2").Value = "2016-08-24 12:50:55"
$wsDeDe.Range("K2").Value = "2016-08-24 12:51:28"

# Overview: Latest HO Xliff Generate Date (G2) -- the max of the handoff
# datetimes across the language sheets for this file.
$wsOverview.Range("G2").Value = "2016-08-24 12:50:55"
